$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 40
$ws.Cells.Item(40, 8).Value = 3092.4167
$ws.Cells.Item(40, 10).Value = 7750
$ws.Cells.Item(40, 12).Value = 7750
$ws.Cells.Item(40, 14).Value = -8100

# ALC row 53
$ws.Cells.Item(53, 8).Value = 448.5
$ws.Cells.Item(53, 9).Value = 499.5
$ws.Cells.Item(53, 10).Value = 423
$ws.Cells.Item(53, 11).Value = 499.5
$ws.Cells.Item(53, 12).Value = 423
$ws.Cells.Item(53, 13).Value = 137.5
$ws.Cells.Item(53, 14).Value = -1697

# ALC row 70
$ws.Cells.Item(70, 8).Value = 1738.9
$ws.Cells.Item(70, 9).Value = 1466.3334
$ws.Cells.Item(70, 11).Value = 4399.0002
$ws.Cells.Item(70, 13).Value = -4129.0002

# ALC row 73
$ws.Cells.Item(73, 8).Value = 1738.9
$ws.Cells.Item(73, 9).Value = 1466.3334
$ws.Cells.Item(73, 11).Value = 4399.0002
$ws.Cells.Item(73, 13).Value = -3463.0002

# ALC row 74
$ws.Cells.Item(74, 8).Value = 8878.200000000001
$ws.Cells.Item(74, 9).Value = 8847.416999999999
$ws.Cells.Item(74, 10).Value = 9001.333000000001
$ws.Cells.Item(74, 11).Value = 8847.416999999999
$ws.Cells.Item(74, 12).Value = 9001.333000000001
$ws.Cells.Item(74, 13).Value = -7911.416999999999
$ws.Cells.Item(74, 14).Value = -10873.333

# ALC row 77
$ws.Cells.Item(77, 8).Value = 8878.200000000001
$ws.Cells.Item(77, 9).Value = 8847.416999999999
$ws.Cells.Item(77, 10).Value = 9001.333000000001
$ws.Cells.Item(77, 11).Value = 44237.085
$ws.Cells.Item(77, 12).Value = 45006.665
$ws.Cells.Item(77, 13).Value = -39557.085
$ws.Cells.Item(77, 14).Value = -54366.665

# ALC row 94
$ws.Cells.Item(94, 8).Value = 1454.6666
$ws.Cells.Item(94, 9).Value = 1058.4
$ws.Cells.Item(94, 11).Value = 1058.4
$ws.Cells.Item(94, 13).Value = -607.4000000000001

# ALC row 112
$ws.Cells.Item(112, 8).Value = 2529.1667
$ws.Cells.Item(112, 9).Value = 2107
$ws.Cells.Item(112, 10).Value = 2567.5454
$ws.Cells.Item(112, 11).Value = 6321
$ws.Cells.Item(112, 12).Value = 7702.6362
$ws.Cells.Item(112, 13).Value = -5213
$ws.Cells.Item(112, 14).Value = -9918.636200000001

# ALC row 115
$ws.Cells.Item(115, 8).Value = 1400
$ws.Cells.Item(115, 9).Value = 1400
$ws.Cells.Item(115, 11).Value = 4200
$ws.Cells.Item(115, 13).Value = -2633

# ALC row 116
$ws.Cells.Item(116, 8).Value = 4498
$ws.Cells.Item(116, 9).Value = 4497.6665
$ws.Cells.Item(116, 11).Value = 4497.6665
$ws.Cells.Item(116, 13).Value = -1055.6665

# ALC row 135
$ws.Cells.Item(135, 8).Value = 392.0625
$ws.Cells.Item(135, 9).Value = 422.85715
$ws.Cells.Item(135, 11).Value = 3805.71435
$ws.Cells.Item(135, 13).Value = -1270.71435

$ws = $wb.Worksheets.Item("ARM")
# ARM row 2
$ws.Cells.Item(2, 8).Value = 1948.8462
$ws.Cells.Item(2, 9).Value = 1948.8462
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 1948.8462
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 14).ClearContents()

# ARM row 32
$ws.Cells.Item(32, 8).Value = 9540.951999999999
$ws.Cells.Item(32, 9).Value = 9540.951999999999
$ws.Cells.Item(32, 11).Value = 9540.951999999999
$ws.Cells.Item(32, 13).Value = -9253.951999999999

# ARM row 74
$ws.Cells.Item(74, 8).Value = 11132.143
$ws.Cells.Item(74, 9).Value = 11132.143
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 11).Value = 11132.143
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 14).ClearContents()

# ARM row 77
$ws.Cells.Item(77, 8).Value = 11132.143
$ws.Cells.Item(77, 9).Value = 11132.143
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 11).Value = 55660.715
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 14).ClearContents()

# ARM row 97
$ws.Cells.Item(97, 8).Value = 2339.6428
$ws.Cells.Item(97, 9).Value = 314.0909
$ws.Cells.Item(97, 11).Value = 314.0909
$ws.Cells.Item(97, 13).Value = 181.9091

# ARM row 116
$ws.Cells.Item(116, 8).Value = 1948.8462
$ws.Cells.Item(116, 9).Value = 1948.8462
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 11).Value = 1948.8462
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 14).ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# BSM row 3
$ws.Cells.Item(3, 8).Value = 1948.8462
$ws.Cells.Item(3, 9).Value = 1948.8462
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 1948.8462
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 14).ClearContents()

# BSM row 35
$ws.Cells.Item(35, 8).Value = 60074
$ws.Cells.Item(35, 10).Value = 60074
$ws.Cells.Item(35, 12).Value = 60074
$ws.Cells.Item(35, 14).Value = -60694

# BSM row 94
$ws.Cells.Item(94, 8).Value = 1651.1428
$ws.Cells.Item(94, 9).Value = 1219.3334
$ws.Cells.Item(94, 10).Value = 1975
$ws.Cells.Item(94, 11).Value = 1219.3334
$ws.Cells.Item(94, 12).Value = 1975
$ws.Cells.Item(94, 13).Value = -768.3334
$ws.Cells.Item(94, 14).Value = -2877

# BSM row 105
$ws.Cells.Item(105, 8).Value = 2967.5
$ws.Cells.Item(105, 9).Value = 3234.75
$ws.Cells.Item(105, 10).Value = 1898.5
$ws.Cells.Item(105, 11).Value = 3234.75
$ws.Cells.Item(105, 12).Value = 1898.5
$ws.Cells.Item(105, 13).Value = -1487.75
$ws.Cells.Item(105, 14).Value = -5392.5

# BSM row 134
$ws.Cells.Item(134, 8).Value = 2219.3333
$ws.Cells.Item(134, 9).Value = 2072.923
$ws.Cells.Item(134, 11).Value = 6218.768999999999
$ws.Cells.Item(134, 13).Value = -3683.768999999999

$ws = $wb.Worksheets.Item("CRP")
# CRP row 58
$ws.Cells.Item(58, 8).Value = 5624
$ws.Cells.Item(58, 9).Value = 2348
$ws.Cells.Item(58, 10).Value = 8900
$ws.Cells.Item(58, 11).Value = 2348
$ws.Cells.Item(58, 12).Value = 8900
$ws.Cells.Item(58, 13).Value = -2145
$ws.Cells.Item(58, 14).Value = -9306

# CRP row 99
$ws.Cells.Item(99, 8).Value = 6999.6
$ws.Cells.Item(99, 9).Value = 7499.75
$ws.Cells.Item(99, 10).Value = 4999
$ws.Cells.Item(99, 11).Value = 7499.75
$ws.Cells.Item(99, 12).Value = 4999
$ws.Cells.Item(99, 13).Value = -6001.75
$ws.Cells.Item(99, 14).Value = -7995

# CRP row 107
$ws.Cells.Item(107, 8).Value = 2064.5
$ws.Cells.Item(107, 9).Value = 2282.6365
$ws.Cells.Item(107, 11).Value = 2282.6365
$ws.Cells.Item(107, 13).Value = -362.6365000000001

# CRP row 126
$ws.Cells.Item(126, 8).Value = 6999.6
$ws.Cells.Item(126, 9).Value = 7499.75
$ws.Cells.Item(126, 10).Value = 4999
$ws.Cells.Item(126, 11).Value = 22499.25
$ws.Cells.Item(126, 12).Value = 14997
$ws.Cells.Item(126, 13).Value = -20029.25
$ws.Cells.Item(126, 14).Value = -19937

# CRP row 136
$ws.Cells.Item(136, 8).Value = 5624
$ws.Cells.Item(136, 9).Value = 2348
$ws.Cells.Item(136, 10).Value = 8900
$ws.Cells.Item(136, 11).Value = 7044
$ws.Cells.Item(136, 12).Value = 26700
$ws.Cells.Item(136, 13).Value = -4494
$ws.Cells.Item(136, 14).Value = -31800

$ws = $wb.Worksheets.Item("CUL")
# CUL row 132
$ws.Cells.Item(132, 8).Value = 959.8333
$ws.Cells.Item(132, 9).Value = 1067.5
$ws.Cells.Item(132, 11).Value = 9607.5
$ws.Cells.Item(132, 13).Value = -7077.5

$ws = $wb.Worksheets.Item("GSM")
# GSM row 102
$ws.Cells.Item(102, 8).Value = 2109.625
$ws.Cells.Item(102, 9).Value = 2109.625
$ws.Cells.Item(102, 11).Value = 2109.625
$ws.Cells.Item(102, 13).Value = -487.625

# GSM row 132
$ws.Cells.Item(132, 8).Value = 3016.5
$ws.Cells.Item(132, 9).Value = 2915
$ws.Cells.Item(132, 11).Value = 8745
$ws.Cells.Item(132, 13).Value = -6215

$ws = $wb.Worksheets.Item("LTW")
# LTW row 16
$ws.Cells.Item(16, 8).Value = 2299
$ws.Cells.Item(16, 9).Value = 2299
$ws.Cells.Item(16, 11).Value = 2299
$ws.Cells.Item(16, 13).Value = -2129

# LTW row 40
$ws.Cells.Item(40, 8).Value = 5098.4287
$ws.Cells.Item(40, 9).Value = 5137.8
$ws.Cells.Item(40, 10).Value = 5000
$ws.Cells.Item(40, 11).Value = 5137.8
$ws.Cells.Item(40, 12).Value = 5000
$ws.Cells.Item(40, 13).Value = -5001.8
$ws.Cells.Item(40, 14).Value = -5272

# LTW row 122
$ws.Cells.Item(122, 8).Value = 6959.8
$ws.Cells.Item(122, 9).Value = 6959.8
$ws.Cells.Item(122, 11).Value = 20879.4
$ws.Cells.Item(122, 13).Value = -18429.4

$ws = $wb.Worksheets.Item("WVR")
# WVR row 122
$ws.Cells.Item(122, 8).Value = 4224.75
$ws.Cells.Item(122, 9).Value = 3966.3333
$ws.Cells.Item(122, 11).Value = 11898.9999
$ws.Cells.Item(122, 13).Value = -9448.999899999999

# WVR row 136
$ws.Cells.Item(136, 8).Value = 6116.8237
$ws.Cells.Item(136, 9).Value = 5231.2144
$ws.Cells.Item(136, 10).Value = 10249.667
$ws.Cells.Item(136, 11).Value = 15693.6432
$ws.Cells.Item(136, 12).Value = 30749.001
$ws.Cells.Item(136, 13).Value = -13143.6432
$ws.Cells.Item(136, 14).Value = -35849.001
